# Applies the "calculation of new indicators" edit described by the diff:
#  - SCHEME_MEASURES: renumber indicator codes MQMSxx -> MQME00x
#  - METADATA_ISSUES: renumber indicator codes MQME10 -> MQME012, MQME01 -> MQME008
#  - METADATA_MEASURES: renumber/renarrate rows and drop the last (now-merged) row
#  - METADATA_METRICS: renumber rows, change some descriptions/values, and append 4 new rows

$wb = $excel.ActiveWorkbook

# Helper: assign a text value to a cell while preventing Excel from
# auto-converting percentage-looking strings ("97.26%", "0.00%", ...) into
# a formatted number. Forcing the Text number format first keeps the
# stored cell type as a string, matching the source data.
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
}

# ---------------------------------------------------------------
# Sheet: SCHEME_MEASURES  (MQMS01..MQMS05 -> MQME001..MQME005)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("SCHEME_MEASURES")
$ws.Range("A2").Value = "MQME001"
$ws.Range("A3").Value = "MQME002"
$ws.Range("A4").Value = "MQME003"
$ws.Range("A5").Value = "MQME004"
$ws.Range("A6").Value = "MQME005"

# ---------------------------------------------------------------
# Sheet: METADATA_ISSUES  (MQME10 -> MQME012, MQME01 -> MQME008)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("METADATA_ISSUES")
$ws.Range("A2").Value = "MQME012"
$ws.Range("A3").Value = "MQME012"
for ($r = 4; $r -le 18; $r++) {
    $ws.Range("A$r").Value = "MQME008"
}

# ---------------------------------------------------------------
# Sheet: METADATA_MEASURES
#   row2: MQME00 / Total number of columns / 362
#      -> MQME006 / Total number of length-required columns / 69
#   row3: MQMEA1 / Total number of length-required columns / 69
#      -> MQME007 / Total number of NUMBER columns / 235
#   row4: MQMEA2 / Total number of NUMBER columns / 235  -> removed (merged into row3)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("METADATA_MEASURES")
$ws.Range("A2").Value = "MQME006"
$ws.Range("B2").Value = "Total number of length-required columns"
$ws.Range("C2").Value = 69
$ws.Range("A3").Value = "MQME007"
$ws.Range("B3").Value = "Total number of NUMBER columns"
$ws.Range("C3").Value = 235
$ws.Rows.Item(4).Delete()

# ---------------------------------------------------------------
# Sheet: METADATA_METRICS
#   IQMEnn -> MQIDnnn, several descriptions change, two new values,
#   and four brand-new rows (9-12) are appended.
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("METADATA_METRICS")

$ws.Range("A2").Value = "MQID001"
$ws.Range("B2").Value = "Table names in singular"
Set-TextValue $ws.Range("C2") "97.26%"

$ws.Range("A3").Value = "MQID002"
$ws.Range("B3").Value = "Table with recommended name length"
# C3 stays "100.00%" (unchanged)

$ws.Range("A4").Value = "MQID003"
$ws.Range("B4").Value = "Columns with correct prefixes"
# C4 stays "100.00%" (unchanged)

$ws.Range("A5").Value = "MQID004"
$ws.Range("B5").Value = "Columns with recommended name size"
# C5 stays "100.00%" (unchanged)

$ws.Range("A6").Value = "MQID005"
$ws.Range("B6").Value = "Columns with comments"
Set-TextValue $ws.Range("C6") "95.86%"

$ws.Range("A7").Value = "MQID006"
$ws.Range("B7").Value = "Table with standard PK prefixes"
# C7 stays "100.00%" (unchanged)

$ws.Range("A8").Value = "MQID007"
$ws.Range("B8").Value = "Table with standard FK prefixes"
# C8 stays "100.00%" (unchanged)

$ws.Range("A9").Value = "MQID008"
$ws.Range("B9").Value = "Table with standard UK prefixes"
Set-TextValue $ws.Range("C9") "0.00%"

$ws.Range("A10").Value = "MQID009"
$ws.Range("B10").Value = "NUMBER columns with valid scale"
Set-TextValue $ws.Range("C10") "100.00%"

$ws.Range("A11").Value = "MQID010"
$ws.Range("B11").Value = "Columns with valid num_distinct"
Set-TextValue $ws.Range("C11") "100.00%"

$ws.Range("A12").Value = "MQID011"
$ws.Range("B12").Value = "Columns with valid num_nulls"
Set-TextValue $ws.Range("C12") "100.00%"
